$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 653.72
$ws.Range("I15").Value = 653.72
$ws.Range("K15").Value = 1961.16
$ws.Range("M15").Value = -1792.16

$ws.Range("H112").Value = 1465
$ws.Range("J112").Value = 1423.125
$ws.Range("L112").Value = 4269.375
$ws.Range("N112").Value = -6485.375

$ws.Range("H133").Value = 88665.664
$ws.Range("J133").Value = 88665.664
$ws.Range("L133").Value = 88665.664
$ws.Range("N133").Value = -98785.664

$ws.Range("H135").Value = 2326.9524
$ws.Range("I135").Value = 181.53334
$ws.Range("J135").Value = 7690.5
$ws.Range("K135").Value = 1633.80006
$ws.Range("L135").Value = 69214.5
$ws.Range("M135").Value = 901.19994
$ws.Range("N135").Value = -74284.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1495.1444
$ws.Range("I32").Value = 827.7308
$ws.Range("K32").Value = 827.7308
$ws.Range("M32").Value = -540.7308

$ws.Range("H97").Value = 1196.1052
$ws.Range("I97").Value = 1201.8572
$ws.Range("J97").Value = 1180
$ws.Range("K97").Value = 1201.8572
$ws.Range("L97").Value = 1180
$ws.Range("M97").Value = -705.8571999999999
$ws.Range("N97").Value = -2172

$ws.Range("H110").Value = 2287.8333
$ws.Range("I110").Value = 2145.4
$ws.Range("K110").Value = 2145.4
$ws.Range("M110").Value = -100.4000000000001

$ws.Range("H132").Value = 3147.2
$ws.Range("I132").Value = 3096.2144
$ws.Range("K132").Value = 9288.643199999999
$ws.Range("M132").Value = -6758.643199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1500
$ws.Range("I86").Value = 1500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -377

$ws.Range("H89").Value = 1500
$ws.Range("I89").Value = 1500
$ws.Range("K89").Value = 7500
$ws.Range("M89").Value = -1884

$ws.Range("H105").Value = 1472.1111
$ws.Range("I105").Value = 1472.1111
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1472.1111
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 274.8888999999999
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 12977.8
$ws.Range("I107").Value = 15599.75
$ws.Range("K107").Value = 15599.75
$ws.Range("M107").Value = -13679.75

$ws.Range("H134").Value = 1825.76
$ws.Range("I134").Value = 1138.9333
$ws.Range("J134").Value = 3886.24
$ws.Range("K134").Value = 3416.7999
$ws.Range("L134").Value = 11658.72
$ws.Range("M134").Value = -881.7999
$ws.Range("N134").Value = -16728.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9378.08
$ws.Range("I31").Value = 13133.125
$ws.Range("J31").Value = 7611
$ws.Range("K31").Value = 13133.125
$ws.Range("L31").Value = 7611
$ws.Range("M31").Value = -12838.125
$ws.Range("N31").Value = -8201

$ws.Range("H34").Value = 9378.08
$ws.Range("I34").Value = 13133.125
$ws.Range("J34").Value = 7611
$ws.Range("K34").Value = 13133.125
$ws.Range("L34").Value = 7611
$ws.Range("M34").Value = -12931.125
$ws.Range("N34").Value = -8015

$ws.Range("H99").Value = 2737.125
$ws.Range("I99").Value = 1999.25
$ws.Range("J99").Value = 3475
$ws.Range("K99").Value = 1999.25
$ws.Range("L99").Value = 3475
$ws.Range("M99").Value = -501.25
$ws.Range("N99").Value = -6471

$ws.Range("H107").Value = 621.5
$ws.Range("I107").Value = 387.25
$ws.Range("K107").Value = 387.25
$ws.Range("M107").Value = 1532.75

$ws.Range("H126").Value = 2737.125
$ws.Range("I126").Value = 1999.25
$ws.Range("J126").Value = 3475
$ws.Range("K126").Value = 5997.75
$ws.Range("L126").Value = 10425
$ws.Range("M126").Value = -3527.75
$ws.Range("N126").Value = -15365

$ws.Range("H127").Value = 59984.715
$ws.Range("J127").Value = 59984.715
$ws.Range("L127").Value = 59984.715
$ws.Range("N127").Value = -69904.715

$ws.Range("H132").Value = 1632.6875
$ws.Range("I132").Value = 1632.6875
$ws.Range("K132").Value = 4898.0625
$ws.Range("M132").Value = -2368.0625

$ws.Range("H141").Value = 84149
$ws.Range("J141").Value = 84149
$ws.Range("L141").Value = 84149
$ws.Range("N141").Value = -94509

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1500
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4040

$ws.Range("H96").Value = 7434.5
$ws.Range("J96").Value = 8021.4
$ws.Range("L96").Value = 24064.2
$ws.Range("N96").Value = -28182.2

$ws.Range("H139").Value = 79769.234
$ws.Range("I139").Value = 93509.63
$ws.Range("K139").Value = 280528.89
$ws.Range("M139").Value = -275388.89

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3533.5
$ws.Range("I107").Value = 5052
$ws.Range("K107").Value = 5052
$ws.Range("M107").Value = -3132

$ws.Range("H123").Value = 59533
$ws.Range("J123").Value = 59799.5
$ws.Range("L123").Value = 59799.5
$ws.Range("N123").Value = -64699.5

$ws.Range("H132").Value = 32270154
$ws.Range("I132").Value = 40008468
$ws.Range("J132").Value = 27182.166
$ws.Range("K132").Value = 120025404
$ws.Range("L132").Value = 81546.49800000001
$ws.Range("M132").Value = -120022874
$ws.Range("N132").Value = -86606.49800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5149.5625
$ws.Range("I7").Value = 2854.182
$ws.Range("K7").Value = 2854.182
$ws.Range("M7").Value = -2742.182

$ws.Range("H68").Value = 2597.762
$ws.Range("I68").Value = 2084
$ws.Range("J68").Value = 4241.8
$ws.Range("K68").Value = 2084
$ws.Range("L68").Value = 4241.8
$ws.Range("M68").Value = -1335
$ws.Range("N68").Value = -5739.8

$ws.Range("H71").Value = 2597.762
$ws.Range("I71").Value = 2084
$ws.Range("J71").Value = 4241.8
$ws.Range("K71").Value = 10420
$ws.Range("L71").Value = 21209
$ws.Range("M71").Value = -6676
$ws.Range("N71").Value = -28697

$ws.Range("H100").Value = 17243.428
$ws.Range("I100").Value = 3450
$ws.Range("K100").Value = 3450
$ws.Range("M100").Value = -2909

$ws.Range("H122").Value = 4496.8
$ws.Range("I122").Value = 3700.28
$ws.Range("J122").Value = 6488.1
$ws.Range("K122").Value = 11100.84
$ws.Range("L122").Value = 19464.3
$ws.Range("M122").Value = -8650.84
$ws.Range("N122").Value = -24364.3

$ws.Range("H126").Value = 5149.5625
$ws.Range("I126").Value = 2854.182
$ws.Range("K126").Value = 8562.545999999998
$ws.Range("M126").Value = -6092.545999999998

$ws.Range("H132").Value = 1905.1111
$ws.Range("I132").Value = 1919.3334
$ws.Range("J132").Value = 1834
$ws.Range("K132").Value = 5758.0002
$ws.Range("L132").Value = 5502
$ws.Range("M132").Value = -3228.0002
$ws.Range("N132").Value = -10562

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 18944.25
$ws.Range("I51").Value = 11925.667
$ws.Range("J51").Value = 40000
$ws.Range("K51").Value = 11925.667
$ws.Range("L51").Value = 40000
$ws.Range("M51").Value = -11415.667
$ws.Range("N51").Value = -41020

$ws.Range("H62").Value = 85108.27
$ws.Range("I62").Value = 182419
$ws.Range("K62").Value = 182419
$ws.Range("M62").Value = -181795

$ws.Range("H65").Value = 85108.27
$ws.Range("I65").Value = 182419
$ws.Range("K65").Value = 912095
$ws.Range("M65").Value = -908975

$ws.Range("H96").Value = 1336
$ws.Range("I96").Value = 899.25
$ws.Range("J96").Value = 1772.75
$ws.Range("K96").Value = 899.25
$ws.Range("L96").Value = 1772.75
$ws.Range("M96").Value = 473.75
$ws.Range("N96").Value = -4518.75

$ws.Range("H104").Value = 396000
$ws.Range("J104").Value = 396000
$ws.Range("L104").Value = 396000
$ws.Range("N104").Value = -402988

$ws.Range("H105").Value = 57391.53
$ws.Range("J105").Value = 57391.53
$ws.Range("L105").Value = 57391.53
$ws.Range("N105").Value = -64379.53

$ws.Range("H122").Value = 1355.6538
$ws.Range("I122").Value = 1356
$ws.Range("K122").Value = 4068
$ws.Range("M122").Value = -1618

$ws.Range("H126").Value = 2678.1765
$ws.Range("I126").Value = 2139.3333
$ws.Range("J126").Value = 3971.4
$ws.Range("K126").Value = 6417.999899999999
$ws.Range("L126").Value = 11914.2
$ws.Range("M126").Value = -3947.999899999999
$ws.Range("N126").Value = -16854.2

$ws.Range("H132").Value = 2685.84
$ws.Range("I132").Value = 1824.3334
$ws.Range("J132").Value = 4901.143
$ws.Range("K132").Value = 5473.0002
$ws.Range("L132").Value = 14703.429
$ws.Range("M132").Value = -2943.0002
$ws.Range("N132").Value = -19763.429

$ws.Range("H136").Value = 3003.2454
$ws.Range("I136").Value = 2421.2708
$ws.Range("K136").Value = 7263.812399999999
$ws.Range("M136").Value = -4713.812399999999

